$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.924.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.64%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.864.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.77%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4358"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3726"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07468"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9346"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.32"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.898.81"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.733"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.444"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06864"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "81.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009058"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.39%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.919.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.119"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.99%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.113.68"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.004"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.51"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.53%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.465"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.714"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09019"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8208"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.817"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.177"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.972"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.54%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("E37").Value = "  -2.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05495"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01979"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.983"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5263"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.038"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1704"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.769"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06754"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.12%  "
$ws.Range("E46").Value = "  -5.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.677"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.881"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -14.80%  "
